$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "29.256.08"
Set-TextValue $ws.Range("E2") "  +0.55%  "
Set-TextValue $ws.Range("D3") "1.858.34"
Set-TextValue $ws.Range("E3") "  +0.46%  "
Set-TextValue $ws.Range("D4") "0.9998"
Set-TextValue $ws.Range("E4") "  +0.17%  "
Set-TextValue $ws.Range("D5") "0.7070"
Set-TextValue $ws.Range("E5") "  +1.75%  "
Set-TextValue $ws.Range("D6") "238.15"
Set-TextValue $ws.Range("D7") "0.9998"
Set-TextValue $ws.Range("E7") "  +0.03%  "
Set-TextValue $ws.Range("D8") "0.08014"
Set-TextValue $ws.Range("E8") "  +3.67%  "
Set-TextValue $ws.Range("D9") "0.3028"
Set-TextValue $ws.Range("E9") "  -0.48%  "
Set-TextValue $ws.Range("D10") "23.53"
Set-TextValue $ws.Range("E10") "  +1.00%  "
Set-TextValue $ws.Range("D11") "0.08184"
Set-TextValue $ws.Range("E11") "  +0.80%  "
Set-TextValue $ws.Range("D12") "1.862.84"
Set-TextValue $ws.Range("E12") "  +0.69%  "
Set-TextValue $ws.Range("D13") "5.196"
Set-TextValue $ws.Range("E13") "  -0.24%  "
Set-TextValue $ws.Range("D14") "0.7064"
Set-TextValue $ws.Range("E14") "  -2.70%  "
Set-TextValue $ws.Range("D15") "89.72"
Set-TextValue $ws.Range("E15") "  +0.80%  "
Set-TextValue $ws.Range("D16") "29.257.15"
Set-TextValue $ws.Range("E16") "  +0.48%  "
Set-TextValue $ws.Range("D17") "0.000007936"
Set-TextValue $ws.Range("E17") "  +1.19%  "
Set-TextValue $ws.Range("D18") "5.819"
Set-TextValue $ws.Range("E18") "  +1.29%  "
Set-TextValue $ws.Range("D19") "13.26"
Set-TextValue $ws.Range("E19") "  +0.44%  "
Set-TextValue $ws.Range("D20") "238.59"
Set-TextValue $ws.Range("E20") "  +1.03%  "
Set-TextValue $ws.Range("D21") "0.9984"
Set-TextValue $ws.Range("E21") "  -0.09%  "
Set-TextValue $ws.Range("D22") "2.107.59"
Set-TextValue $ws.Range("E22") "  +0.14%  "
Set-TextValue $ws.Range("D23") "0.9997"
Set-TextValue $ws.Range("E23") "  +0.04%  "
Set-TextValue $ws.Range("D24") "7.478"
Set-TextValue $ws.Range("D25") "162.89"
Set-TextValue $ws.Range("E25") "  +1.14%  "
Set-TextValue $ws.Range("D26") "8.893"
Set-TextValue $ws.Range("E26") "  -0.99%  "
Set-TextValue $ws.Range("D27") "0.1437"
Set-TextValue $ws.Range("E27") "  +0.35%  "
Set-TextValue $ws.Range("D28") "18.12"
Set-TextValue $ws.Range("E28") "  +0.32%  "
Set-TextValue $ws.Range("D29") "1.928"
Set-TextValue $ws.Range("E29") "  -2.36%  "
Set-TextValue $ws.Range("D30") "1.427"
Set-TextValue $ws.Range("E30") "  +2.28%  "
Set-TextValue $ws.Range("E31") "  -0.58%  "
Set-TextValue $ws.Range("D32") "4.374"
Set-TextValue $ws.Range("E32") "  -2.73%  "
Set-TextValue $ws.Range("D33") "4.026"
Set-TextValue $ws.Range("E33") "  +0.40%  "
Set-TextValue $ws.Range("E34") "  -0.73%  "
Set-TextValue $ws.Range("E35") "  -1.80%  "
Set-TextValue $ws.Range("D36") "0.7151"
Set-TextValue $ws.Range("E36") "  +1.37%  "
Set-TextValue $ws.Range("D37") "1.002"
Set-TextValue $ws.Range("E37") "  -1.91%  "
Set-TextValue $ws.Range("D38") "2.661"
Set-TextValue $ws.Range("E38") "  +0.63%  "
Set-TextValue $ws.Range("D39") "0.01857"
Set-TextValue $ws.Range("E39") "  +0.07%  "
Set-TextValue $ws.Range("D40") "2.724"
Set-TextValue $ws.Range("D41") "0.9370"
Set-TextValue $ws.Range("E41") "  +2.74%  "
Set-TextValue $ws.Range("D42") "1.139.15"
Set-TextValue $ws.Range("E42") "  +3.93%  "
Set-TextValue $ws.Range("D43") "5.983"
Set-TextValue $ws.Range("E43") "  -0.51%  "
Set-TextValue $ws.Range("D44") "0.4270"
Set-TextValue $ws.Range("E44") "  -0.06%  "
Set-TextValue $ws.Range("D45") "70.57"
Set-TextValue $ws.Range("E45") "  -0.10%  "
Set-TextValue $ws.Range("D46") "0.9994"
Set-TextValue $ws.Range("E46") "  +0.11%  "
Set-TextValue $ws.Range("D47") "102.75"
Set-TextValue $ws.Range("E47") "  -0.20%  "
Set-TextValue $ws.Range("D48") "0.5303"
Set-TextValue $ws.Range("E48") "  -4.24%  "
Set-TextValue $ws.Range("D49") "1.764"
Set-TextValue $ws.Range("E49") "  -0.31%  "
Set-TextValue $ws.Range("D50") "1.997.38"
Set-TextValue $ws.Range("E50") "  +0.00%  "
Set-TextValue $ws.Range("D51") "9.182"
Set-TextValue $ws.Range("E51") "  +0.28%  "
